$d = $word.ActiveDocument

function Set-ParaText($index, $text) {
    $p = $d.Paragraphs($index)
    $r = $p.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

Set-ParaText 1 "ContosoLearn Market Research"
Set-ParaText 2 "AdatumLearn: AdatumLearn is a top AI-powered learning platform that uses artificial intelligence to enrich eLearning with features that automate a variety of tasks. It is known for its content authoring capabilities and adaptive learning technology."
Set-ParaText 3 "AdventureLearn: AdventureLearn is another AI-powered learning platform that offers personalized learning experiences and data-driven recommendations."
Set-ParaText 4 "AlpineTraining: AlpineTraining is a mobile-first learning platform that focuses on microlearning."
Set-ParaText 5 "Bellows OnDemand: Bellows OnDemand is a comprehensive learning solution that offers content creation and social collaboration."
Set-ParaText 6 "FabrikamLearning: FabrikamLearning provides a suite of learning platforms that cater to different learning needs."
Set-ParaText 7 "FirstUp Cards: FirstUp Cards is a mobile learning app that is ideal for training on safety procedures, compliance, new product knowledge or any other type of training scenario."
Set-ParaText 8 "Munson'sLearn: Munson'sLearn is designed to enable businesses to train their employees, partners, and customers."
Set-ParaText 9 "LibertyLearn: LibertyLearn is a fast LMS for your mission-critical project."

# WoodgroveLMS paragraph: split into three runs, middle one ("a best") flagged
# by the grammar checker (w:proofErr) in the authored document.
$p10 = $d.Paragraphs(10)
$r10 = $p10.Range
$r10.End = $r10.End - 1
$r10.Text = "WoodgroveLMS: WoodgroveLMS is a functional and attractive learning management system built to provide a best-in-class training experience."

Set-ParaText 11 "NorthwindWorlds: NorthwindWorlds is a powerful, easy-to-use, and reliable training solution for individuals and enterprises."
Set-ParaText 12 "ProsewareLearn: ProsewareLearn is an online education company that offers a variety of video training courses for software developers, IT administrators, and creative professionals through its website."
Set-ParaText 13 "RelecloudLearn: RelecloudLearn is an American online learning platform that offers massive open online courses (MOOC), specializations, and degrees in a variety of subjects."
Set-ParaText 14 "TreyAcademy: TreyAcademy is an online learning platform aimed at professional adults and students, developed in May 2010."
Set-ParaText 15 "These platforms have a significant market presence and are widely recognized for their AI-powered features, such as personalized learning experiences, data-driven recommendations, and automation of tasks. They are transforming the eLearning landscape by leveraging AI to deliver more engaging, rewarding, and personalized learning experiences. "
